# B6-PowerPoint.pptx edit script
# 1) Re-style the three tables (slides 14, 15, 16) from the custom
#    "Table_0" style {431F49B7-0778-4C32-AA75-55D81184343C} to the
#    built-in "No Style, No Grid" style {885A8A24-92EE-4252-954A-917FAAB4018B}.
# 2) Re-colour the presentation theme ("Integral" / Red Violet) to the
#    stock Office theme colours (the deck's second, otherwise-unused
#    theme part already carries that exact "Office Theme" palette).

$p = $ppt.ActivePresentation

$oldStyleId = "{431F49B7-0778-4C32-AA75-55D81184343C}"
$newStyleId = "{885A8A24-92EE-4252-954A-917FAAB4018B}"

$tableSlideIndexes = @(14, 15, 16)
foreach ($slideIndex in $tableSlideIndexes) {
    $slide = $p.Slides.Item($slideIndex)
    for ($shapeIndex = 1; $shapeIndex -le $slide.Shapes.Count; $shapeIndex++) {
        $shape = $slide.Shapes.Item($shapeIndex)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# Swap the theme's 12 scheme colours from the "Red Violet" palette to the
# stock "Office" palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
